# Automatic update of files.
# Applies the per-row data refresh to rows 2-8 of the Artfynd sheet as
# described by the source diff. Each row's species/observation data is
# rewritten in place (row positions stay fixed; contents are refreshed),
# and a new "Aktivitet" (M4) value is added for row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 111739317
$ws.Range("B2").Value = 78579
$ws.Range("E2").Value = 2081
$ws.Range("F2").Value = "Skrovellav"
$ws.Range("G2").Value = "Lobaria scrobiculata"
$ws.Range("H2").Value = "(Scop.) DC."
$ws.Range("Q2").Value = 573911.5177193542
$ws.Range("R2").Value = 7172648.020174325

# --- Row 3 ---
$ws.Range("A3").Value = 111739316
$ws.Range("B3").Value = 78578
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = "Lunglav"
$ws.Range("G3").Value = "Lobaria pulmonaria"
$ws.Range("H3").Value = "(L.) Hoffm."
$ws.Range("Q3").Value = 573904.5013778479
$ws.Range("R3").Value = 7172636.708955797

# --- Row 4 ---
$ws.Range("A4").Value = 111739306
$ws.Range("B4").Value = 56398
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 573906.0397215446
$ws.Range("R4").Value = 7172521.061635921

# --- Row 5 ---
$ws.Range("A5").Value = 111739309
$ws.Range("B5").Value = 78536
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 229497
$ws.Range("F5").Value = "Korallblylav"
$ws.Range("G5").Value = "Parmeliella triptophylla"
$ws.Range("H5").Value = "(Ach.) Müll.Arg."
$ws.Range("Q5").Value = 574011.1276117128
$ws.Range("R5").Value = 7172434.078971106

# --- Row 6 ---
$ws.Range("A6").Value = 111739315
$ws.Range("B6").Value = 78605
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 6462
$ws.Range("F6").Value = "Stuplav"
$ws.Range("G6").Value = "Nephroma bellum"
$ws.Range("H6").Value = "(Spreng.) Tuck."
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 573904.5013778479
$ws.Range("R6").Value = 7172636.708955797

# --- Row 7 ---
$ws.Range("A7").Value = 111739311
$ws.Range("B7").Value = 77515
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 574011.8892867711
$ws.Range("R7").Value = 7172473.089384713

# --- Row 8 ---
$ws.Range("A8").Value = 111739313
$ws.Range("B8").Value = 73701
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1467
$ws.Range("F8").Value = "Rödbrun blekspik"
$ws.Range("G8").Value = "Sclerophora coniophaea"
$ws.Range("H8").Value = "(Norman) J.Mattsson & Middelb."
$ws.Range("Q8").Value = 574025.0565134182
$ws.Range("R8").Value = 7172443.417908707
